# Updated symbol list on Sun Dec 11 23:50:19 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numeric-looking values as TEXT in the
# original workbook. Assigning a plain numeric-looking string to
# Range.Value lets Excel auto-convert it to a real number, which would
# change the cell's stored type. To keep it text (matching the source
# data) without leaving a stray "quote-prefixed" style behind, we
# temporarily force a text number format, assign the value, then put
# the style back to Normal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "284.83"
Set-TextValue "D3"  "21.10"
Set-TextValue "D4"  "6.446"
Set-TextValue "D5"  "0.06346"
Set-TextValue "D6"  "3.602"
Set-TextValue "D7"  "1.526"
Set-TextValue "D8"  "6.557"
Set-TextValue "D9"  "0.8202"
Set-TextValue "D12" "0.08620"
Set-TextValue "D13" "0.03668"
Set-TextValue "D14" "0.03217"
Set-TextValue "D16" "3.716"
Set-TextValue "D17" "0.001641"
Set-TextValue "D18" "0.04734"
Set-TextValue "D19" "0.006131"
Set-TextValue "D20" "0.006272"
Set-TextValue "D22" "0.0001602"
Set-TextValue "D23" "3.783"
Set-TextValue "D25" "0.3354"
Set-TextValue "D26" "0.1262"
Set-TextValue "D40" "0.04759"
Set-TextValue "D41" "0.007130"

# --- Rows 42 / 43: BKEXToken and CEJI swapped places in the ranking ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.004506"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1113"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# --- Remaining price / label updates ---
Set-TextValue "D44" "0.01145"
Set-TextValue "D45" "0.00006739"
Set-TextValue "D47" "1.001"

Set-TextValue "D48" "0.003181"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

Set-TextValue "D49" "0.00001502"
$ws.Range("E49").Value = "48CryptobidCoinCBC"
